$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shape = $s.Shapes.Item(2)

$lines = @(
    "Path following",
    "calculate distance to goal  ",
    "if distance is above threshold:  ",
    "    calculate angle to goal",
    "    if angle is above threshold:",
    "        re-orient robot",
    "    else:",
    "        move robot",
    "else:",
    "    increment goal index in path array",
    "    if index surpasses the length of the path:",
    "        stop the movement timer",
    ""
)

$text = [string]::Join([char]13, $lines)
$shape.TextFrame.TextRange.Text = $text
